# Commit message: "add the NA's under duplicate_image_filename"
#
# Column E ("duplicate_image_filename", header already present in E1) was
# empty for the data rows. This fills rows 2-21 (the practice + main trial
# rows) with the literal string "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2:E21").Value = "NA"
